$d = $word.ActiveDocument

# Update the date/title line
$d.Paragraphs.Item(1).Range.Text = "2024-10-09 Wednesday"

# Update the 20x5 table of arithmetic expressions (positional, row-major)
$t = $d.Tables.Item(1)
$values = @(
    "19+77=96", "81-67=14", "2+2=4", "27+32=59", "44+2=46",
    "89-79=10", "62+36=98", "69-17=52", "65-65=0", "96-89=7",
    "98-90=8", "1+87=88", "74-57=17", "51-8=43", "80-59=21",
    "82-45=37", "48+16=64", "37+22=59", "32+44=76", "27+38=65",
    "24+25=49", "98-78=20", "56-5=51", "26+51=77", "48+20=68",
    "86-7=79", "91-67=24", "44+54=98", "64-27=37", "46+18=64",
    "51+5=56", "40-26=14", "26+50=76", "0+60=60", "44-28=16",
    "55-27=28", "39-19=20", "76+14=90", "78-53=25", "56-17=39",
    "66-56=10", "62-18=44", "31+17=48", "52-6=46", "79-29=50",
    "75+10=85", "68-37=31", "42+13=55", "62-23=39", "40-24=16",
    "17+17=34", "5+21=26", "54+0=54", "22-4=18", "91-7=84",
    "26-3=23", "61-6=55", "28+22=50", "68-62=6", "54-2=52",
    "56+36=92", "82-36=46", "35-3=32", "13+33=46", "83-13=70",
    "62+26=88", "37+51=88", "99-63=36", "26-8=18", "42+44=86",
    "97-36=61", "46+48=94", "96-93=3", "50+7=57", "10+76=86",
    "36-29=7", "78-30=48", "24+61=85", "20+49=69", "16+24=40",
    "65+27=92", "1+4=5", "24+59=83", "49+43=92", "94-35=59",
    "0+96=96", "63+0=63", "35-19=16", "83-74=9", "84-17=67",
    "10+67=77", "81-47=34", "63+22=85", "42+53=95", "17+65=82",
    "4+60=64", "39-11=28", "88-12=76", "82-0=82", "28+7=35"
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}
